# This script re-applies the weekly Fruta/Hortaliza data refresh described in the
# commit "Fruta / hortaliza, semanal": the values of columns D,L,M,N,O,P,Q,R,S,T
# are redistributed across the data rows (rows 2-41). We read every source cell
# with Value2 first (so row data is captured before anything is overwritten), and
# then write the captured values into their destination rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Phase 1: capture current values from their source rows ----
$v_D2 = $ws.Range("D26").Value2
$v_L2 = $ws.Range("L26").Value2
$v_M2 = $ws.Range("M26").Value2
$v_N2 = $ws.Range("N26").Value2
$v_O2 = $ws.Range("O26").Value2
$v_P2 = $ws.Range("P26").Value2
$v_Q2 = $ws.Range("Q26").Value2
$v_R2 = $ws.Range("R26").Value2
$v_S2 = $ws.Range("S26").Value2
$v_T2 = $ws.Range("T26").Value2

$v_D3 = $ws.Range("D3").Value2
$v_L3 = $ws.Range("L3").Value2
$v_M3 = $ws.Range("M3").Value2
$v_N3 = $ws.Range("N3").Value2
$v_O3 = $ws.Range("O3").Value2
$v_P3 = $ws.Range("P3").Value2
$v_Q3 = $ws.Range("Q3").Value2
$v_R3 = $ws.Range("R3").Value2
$v_S3 = $ws.Range("S3").Value2
$v_T3 = $ws.Range("T3").Value2

$v_D4 = $ws.Range("D7").Value2
$v_L4 = $ws.Range("L7").Value2
$v_M4 = $ws.Range("M7").Value2
$v_N4 = $ws.Range("N7").Value2
$v_O4 = $ws.Range("O7").Value2
$v_P4 = $ws.Range("P7").Value2
$v_Q4 = $ws.Range("Q7").Value2
$v_R4 = $ws.Range("R7").Value2
$v_S4 = $ws.Range("S7").Value2
$v_T4 = $ws.Range("T7").Value2

$v_D5 = $ws.Range("D12").Value2
$v_L5 = $ws.Range("L12").Value2
$v_M5 = $ws.Range("M12").Value2
$v_N5 = $ws.Range("N12").Value2
$v_O5 = $ws.Range("O12").Value2
$v_P5 = $ws.Range("P12").Value2
$v_Q5 = $ws.Range("Q12").Value2
$v_R5 = $ws.Range("R12").Value2
$v_S5 = $ws.Range("S12").Value2
$v_T5 = $ws.Range("T12").Value2

$v_D6 = $ws.Range("D27").Value2
$v_L6 = $ws.Range("L27").Value2
$v_M6 = $ws.Range("M27").Value2
$v_N6 = $ws.Range("N27").Value2
$v_O6 = $ws.Range("O27").Value2
$v_P6 = $ws.Range("P27").Value2
$v_Q6 = $ws.Range("Q27").Value2
$v_R6 = $ws.Range("R27").Value2
$v_S6 = $ws.Range("S27").Value2
$v_T6 = $ws.Range("T27").Value2

$v_D7 = $ws.Range("D35").Value2
$v_L7 = $ws.Range("L35").Value2
$v_M7 = $ws.Range("M35").Value2
$v_N7 = $ws.Range("N35").Value2
$v_O7 = $ws.Range("O35").Value2
$v_P7 = $ws.Range("P35").Value2
$v_Q7 = $ws.Range("Q35").Value2
$v_R7 = $ws.Range("R35").Value2
$v_S7 = $ws.Range("S35").Value2
$v_T7 = $ws.Range("T35").Value2

$v_D8 = $ws.Range("D5").Value2
$v_L8 = $ws.Range("L5").Value2
$v_M8 = $ws.Range("M5").Value2
$v_N8 = $ws.Range("N5").Value2
$v_O8 = $ws.Range("O5").Value2
$v_P8 = $ws.Range("P5").Value2
$v_Q8 = $ws.Range("Q5").Value2
$v_R8 = $ws.Range("R5").Value2
$v_S8 = $ws.Range("S5").Value2
$v_T8 = $ws.Range("T5").Value2

$v_D9 = $ws.Range("D6").Value2
$v_L9 = $ws.Range("L6").Value2
$v_M9 = $ws.Range("M6").Value2
$v_N9 = $ws.Range("N6").Value2
$v_O9 = $ws.Range("O6").Value2
$v_P9 = $ws.Range("P6").Value2
$v_Q9 = $ws.Range("Q6").Value2
$v_R9 = $ws.Range("R6").Value2
$v_S9 = $ws.Range("S6").Value2
$v_T9 = $ws.Range("T6").Value2

$v_D10 = $ws.Range("D37").Value2
$v_L10 = $ws.Range("L37").Value2
$v_M10 = $ws.Range("M37").Value2
$v_N10 = $ws.Range("N37").Value2
$v_O10 = $ws.Range("O37").Value2
$v_P10 = $ws.Range("P37").Value2
$v_Q10 = $ws.Range("Q37").Value2
$v_R10 = $ws.Range("R37").Value2
$v_S10 = $ws.Range("S37").Value2
$v_T10 = $ws.Range("T37").Value2

$v_D11 = $ws.Range("D30").Value2
$v_L11 = $ws.Range("L30").Value2
$v_M11 = $ws.Range("M30").Value2
$v_N11 = $ws.Range("N30").Value2
$v_O11 = $ws.Range("O30").Value2
$v_P11 = $ws.Range("P30").Value2
$v_Q11 = $ws.Range("Q30").Value2
$v_R11 = $ws.Range("R30").Value2
$v_S11 = $ws.Range("S30").Value2
$v_T11 = $ws.Range("T30").Value2

$v_D12 = $ws.Range("D36").Value2
$v_L12 = $ws.Range("L36").Value2
$v_M12 = $ws.Range("M36").Value2
$v_N12 = $ws.Range("N36").Value2
$v_O12 = $ws.Range("O36").Value2
$v_P12 = $ws.Range("P36").Value2
$v_Q12 = $ws.Range("Q36").Value2
$v_R12 = $ws.Range("R36").Value2
$v_S12 = $ws.Range("S36").Value2
$v_T12 = $ws.Range("T36").Value2

$v_D13 = $ws.Range("D20").Value2
$v_L13 = $ws.Range("L20").Value2
$v_M13 = $ws.Range("M20").Value2
$v_N13 = $ws.Range("N20").Value2
$v_O13 = $ws.Range("O20").Value2
$v_P13 = $ws.Range("P20").Value2
$v_Q13 = $ws.Range("Q20").Value2
$v_R13 = $ws.Range("R20").Value2
$v_S13 = $ws.Range("S20").Value2
$v_T13 = $ws.Range("T20").Value2

$v_D14 = $ws.Range("D28").Value2
$v_L14 = $ws.Range("L28").Value2
$v_M14 = $ws.Range("M28").Value2
$v_N14 = $ws.Range("N28").Value2
$v_O14 = $ws.Range("O28").Value2
$v_P14 = $ws.Range("P28").Value2
$v_Q14 = $ws.Range("Q28").Value2
$v_R14 = $ws.Range("R28").Value2
$v_S14 = $ws.Range("S28").Value2
$v_T14 = $ws.Range("T28").Value2

$v_D15 = $ws.Range("D18").Value2
$v_L15 = $ws.Range("L18").Value2
$v_M15 = $ws.Range("M18").Value2
$v_N15 = $ws.Range("N18").Value2
$v_O15 = $ws.Range("O18").Value2
$v_P15 = $ws.Range("P18").Value2
$v_Q15 = $ws.Range("Q18").Value2
$v_R15 = $ws.Range("R18").Value2
$v_S15 = $ws.Range("S18").Value2
$v_T15 = $ws.Range("T18").Value2

$v_D16 = $ws.Range("D19").Value2
$v_L16 = $ws.Range("L19").Value2
$v_M16 = $ws.Range("M19").Value2
$v_N16 = $ws.Range("N19").Value2
$v_O16 = $ws.Range("O19").Value2
$v_P16 = $ws.Range("P19").Value2
$v_Q16 = $ws.Range("Q19").Value2
$v_R16 = $ws.Range("R19").Value2
$v_S16 = $ws.Range("S19").Value2
$v_T16 = $ws.Range("T19").Value2

$v_D17 = $ws.Range("D23").Value2
$v_L17 = $ws.Range("L23").Value2
$v_M17 = $ws.Range("M23").Value2
$v_N17 = $ws.Range("N23").Value2
$v_O17 = $ws.Range("O23").Value2
$v_P17 = $ws.Range("P23").Value2
$v_Q17 = $ws.Range("Q23").Value2
$v_R17 = $ws.Range("R23").Value2
$v_S17 = $ws.Range("S23").Value2
$v_T17 = $ws.Range("T23").Value2

$v_D18 = $ws.Range("D10").Value2
$v_L18 = $ws.Range("L10").Value2
$v_M18 = $ws.Range("M10").Value2
$v_N18 = $ws.Range("N10").Value2
$v_O18 = $ws.Range("O10").Value2
$v_P18 = $ws.Range("P10").Value2
$v_Q18 = $ws.Range("Q10").Value2
$v_R18 = $ws.Range("R10").Value2
$v_S18 = $ws.Range("S10").Value2
$v_T18 = $ws.Range("T10").Value2

$v_D19 = $ws.Range("D14").Value2
$v_L19 = $ws.Range("L14").Value2
$v_M19 = $ws.Range("M14").Value2
$v_N19 = $ws.Range("N14").Value2
$v_O19 = $ws.Range("O14").Value2
$v_P19 = $ws.Range("P14").Value2
$v_Q19 = $ws.Range("Q14").Value2
$v_R19 = $ws.Range("R14").Value2
$v_S19 = $ws.Range("S14").Value2
$v_T19 = $ws.Range("T14").Value2

$v_D20 = $ws.Range("D15").Value2
$v_L20 = $ws.Range("L15").Value2
$v_M20 = $ws.Range("M15").Value2
$v_N20 = $ws.Range("N15").Value2
$v_O20 = $ws.Range("O15").Value2
$v_P20 = $ws.Range("P15").Value2
$v_Q20 = $ws.Range("Q15").Value2
$v_R20 = $ws.Range("R15").Value2
$v_S20 = $ws.Range("S15").Value2
$v_T20 = $ws.Range("T15").Value2

$v_D21 = $ws.Range("D40").Value2
$v_L21 = $ws.Range("L40").Value2
$v_M21 = $ws.Range("M40").Value2
$v_N21 = $ws.Range("N40").Value2
$v_O21 = $ws.Range("O40").Value2
$v_P21 = $ws.Range("P40").Value2
$v_Q21 = $ws.Range("Q40").Value2
$v_R21 = $ws.Range("R40").Value2
$v_S21 = $ws.Range("S40").Value2
$v_T21 = $ws.Range("T40").Value2

$v_D22 = $ws.Range("D41").Value2
$v_L22 = $ws.Range("L41").Value2
$v_M22 = $ws.Range("M41").Value2
$v_N22 = $ws.Range("N41").Value2
$v_O22 = $ws.Range("O41").Value2
$v_P22 = $ws.Range("P41").Value2
$v_Q22 = $ws.Range("Q41").Value2
$v_R22 = $ws.Range("R41").Value2
$v_S22 = $ws.Range("S41").Value2
$v_T22 = $ws.Range("T41").Value2

$v_D23 = $ws.Range("D22").Value2
$v_L23 = $ws.Range("L22").Value2
$v_M23 = $ws.Range("M22").Value2
$v_N23 = $ws.Range("N22").Value2
$v_O23 = $ws.Range("O22").Value2
$v_P23 = $ws.Range("P22").Value2
$v_Q23 = $ws.Range("Q22").Value2
$v_R23 = $ws.Range("R22").Value2
$v_S23 = $ws.Range("S22").Value2
$v_T23 = $ws.Range("T22").Value2

$v_D24 = $ws.Range("D16").Value2
$v_L24 = $ws.Range("L16").Value2
$v_M24 = $ws.Range("M16").Value2
$v_N24 = $ws.Range("N16").Value2
$v_O24 = $ws.Range("O16").Value2
$v_P24 = $ws.Range("P16").Value2
$v_Q24 = $ws.Range("Q16").Value2
$v_R24 = $ws.Range("R16").Value2
$v_S24 = $ws.Range("S16").Value2
$v_T24 = $ws.Range("T16").Value2

$v_D25 = $ws.Range("D31").Value2
$v_L25 = $ws.Range("L31").Value2
$v_M25 = $ws.Range("M31").Value2
$v_N25 = $ws.Range("N31").Value2
$v_O25 = $ws.Range("O31").Value2
$v_P25 = $ws.Range("P31").Value2
$v_Q25 = $ws.Range("Q31").Value2
$v_R25 = $ws.Range("R31").Value2
$v_S25 = $ws.Range("S31").Value2
$v_T25 = $ws.Range("T31").Value2

$v_D26 = $ws.Range("D32").Value2
$v_L26 = $ws.Range("L32").Value2
$v_M26 = $ws.Range("M32").Value2
$v_N26 = $ws.Range("N32").Value2
$v_O26 = $ws.Range("O32").Value2
$v_P26 = $ws.Range("P32").Value2
$v_Q26 = $ws.Range("Q32").Value2
$v_R26 = $ws.Range("R32").Value2
$v_S26 = $ws.Range("S32").Value2
$v_T26 = $ws.Range("T32").Value2

$v_D27 = $ws.Range("D25").Value2
$v_L27 = $ws.Range("L25").Value2
$v_M27 = $ws.Range("M25").Value2
$v_N27 = $ws.Range("N25").Value2
$v_O27 = $ws.Range("O25").Value2
$v_P27 = $ws.Range("P25").Value2
$v_Q27 = $ws.Range("Q25").Value2
$v_R27 = $ws.Range("R25").Value2
$v_S27 = $ws.Range("S25").Value2
$v_T27 = $ws.Range("T25").Value2

$v_D28 = $ws.Range("D29").Value2
$v_L28 = $ws.Range("L29").Value2
$v_M28 = $ws.Range("M29").Value2
$v_N28 = $ws.Range("N29").Value2
$v_O28 = $ws.Range("O29").Value2
$v_P28 = $ws.Range("P29").Value2
$v_Q28 = $ws.Range("Q29").Value2
$v_R28 = $ws.Range("R29").Value2
$v_S28 = $ws.Range("S29").Value2
$v_T28 = $ws.Range("T29").Value2

$v_D29 = $ws.Range("D24").Value2
$v_L29 = $ws.Range("L24").Value2
$v_M29 = $ws.Range("M24").Value2
$v_N29 = $ws.Range("N24").Value2
$v_O29 = $ws.Range("O24").Value2
$v_P29 = $ws.Range("P24").Value2
$v_Q29 = $ws.Range("Q24").Value2
$v_R29 = $ws.Range("R24").Value2
$v_S29 = $ws.Range("S24").Value2
$v_T29 = $ws.Range("T24").Value2

$v_D30 = $ws.Range("D17").Value2
$v_L30 = $ws.Range("L17").Value2
$v_M30 = $ws.Range("M17").Value2
$v_N30 = $ws.Range("N17").Value2
$v_O30 = $ws.Range("O17").Value2
$v_P30 = $ws.Range("P17").Value2
$v_Q30 = $ws.Range("Q17").Value2
$v_R30 = $ws.Range("R17").Value2
$v_S30 = $ws.Range("S17").Value2
$v_T30 = $ws.Range("T17").Value2

$v_D31 = $ws.Range("D13").Value2
$v_L31 = $ws.Range("L13").Value2
$v_M31 = $ws.Range("M13").Value2
$v_N31 = $ws.Range("N13").Value2
$v_O31 = $ws.Range("O13").Value2
$v_P31 = $ws.Range("P13").Value2
$v_Q31 = $ws.Range("Q13").Value2
$v_R31 = $ws.Range("R13").Value2
$v_S31 = $ws.Range("S13").Value2
$v_T31 = $ws.Range("T13").Value2

$v_D32 = $ws.Range("D2").Value2
$v_L32 = $ws.Range("L2").Value2
$v_M32 = $ws.Range("M2").Value2
$v_N32 = $ws.Range("N2").Value2
$v_O32 = $ws.Range("O2").Value2
$v_P32 = $ws.Range("P2").Value2
$v_Q32 = $ws.Range("Q2").Value2
$v_R32 = $ws.Range("R2").Value2
$v_S32 = $ws.Range("S2").Value2
$v_T32 = $ws.Range("T2").Value2

$v_D33 = $ws.Range("D34").Value2
$v_L33 = $ws.Range("L34").Value2
$v_M33 = $ws.Range("M34").Value2
$v_N33 = $ws.Range("N34").Value2
$v_O33 = $ws.Range("O34").Value2
$v_P33 = $ws.Range("P34").Value2
$v_Q33 = $ws.Range("Q34").Value2
$v_R33 = $ws.Range("R34").Value2
$v_S33 = $ws.Range("S34").Value2
$v_T33 = $ws.Range("T34").Value2

$v_D34 = $ws.Range("D21").Value2
$v_L34 = $ws.Range("L21").Value2
$v_M34 = $ws.Range("M21").Value2
$v_N34 = $ws.Range("N21").Value2
$v_O34 = $ws.Range("O21").Value2
$v_P34 = $ws.Range("P21").Value2
$v_Q34 = $ws.Range("Q21").Value2
$v_R34 = $ws.Range("R21").Value2
$v_S34 = $ws.Range("S21").Value2
$v_T34 = $ws.Range("T21").Value2

$v_D35 = $ws.Range("D4").Value2
$v_L35 = $ws.Range("L4").Value2
$v_M35 = $ws.Range("M4").Value2
$v_N35 = $ws.Range("N4").Value2
$v_O35 = $ws.Range("O4").Value2
$v_P35 = $ws.Range("P4").Value2
$v_Q35 = $ws.Range("Q4").Value2
$v_R35 = $ws.Range("R4").Value2
$v_S35 = $ws.Range("S4").Value2
$v_T35 = $ws.Range("T4").Value2

$v_D36 = $ws.Range("D8").Value2
$v_L36 = $ws.Range("L8").Value2
$v_M36 = $ws.Range("M8").Value2
$v_N36 = $ws.Range("N8").Value2
$v_O36 = $ws.Range("O8").Value2
$v_P36 = $ws.Range("P8").Value2
$v_Q36 = $ws.Range("Q8").Value2
$v_R36 = $ws.Range("R8").Value2
$v_S36 = $ws.Range("S8").Value2
$v_T36 = $ws.Range("T8").Value2

$v_D37 = $ws.Range("D9").Value2
$v_L37 = $ws.Range("L9").Value2
$v_M37 = $ws.Range("M9").Value2
$v_N37 = $ws.Range("N9").Value2
$v_O37 = $ws.Range("O9").Value2
$v_P37 = $ws.Range("P9").Value2
$v_Q37 = $ws.Range("Q9").Value2
$v_R37 = $ws.Range("R9").Value2
$v_S37 = $ws.Range("S9").Value2
$v_T37 = $ws.Range("T9").Value2

$v_D38 = $ws.Range("D11").Value2
$v_L38 = $ws.Range("L11").Value2
$v_M38 = $ws.Range("M11").Value2
$v_N38 = $ws.Range("N11").Value2
$v_O38 = $ws.Range("O11").Value2
$v_P38 = $ws.Range("P11").Value2
$v_Q38 = $ws.Range("Q11").Value2
$v_R38 = $ws.Range("R11").Value2
$v_S38 = $ws.Range("S11").Value2
$v_T38 = $ws.Range("T11").Value2

$v_D39 = $ws.Range("D39").Value2
$v_L39 = $ws.Range("L39").Value2
$v_M39 = $ws.Range("M39").Value2
$v_N39 = $ws.Range("N39").Value2
$v_O39 = $ws.Range("O39").Value2
$v_P39 = $ws.Range("P39").Value2
$v_Q39 = $ws.Range("Q39").Value2
$v_R39 = $ws.Range("R39").Value2
$v_S39 = $ws.Range("S39").Value2
$v_T39 = $ws.Range("T39").Value2

$v_D40 = $ws.Range("D38").Value2
$v_L40 = $ws.Range("L38").Value2
$v_M40 = $ws.Range("M38").Value2
$v_N40 = $ws.Range("N38").Value2
$v_O40 = $ws.Range("O38").Value2
$v_P40 = $ws.Range("P38").Value2
$v_Q40 = $ws.Range("Q38").Value2
$v_R40 = $ws.Range("R38").Value2
$v_S40 = $ws.Range("S38").Value2
$v_T40 = $ws.Range("T38").Value2

$v_D41 = $ws.Range("D33").Value2
$v_L41 = $ws.Range("L33").Value2
$v_M41 = $ws.Range("M33").Value2
$v_N41 = $ws.Range("N33").Value2
$v_O41 = $ws.Range("O33").Value2
$v_P41 = $ws.Range("P33").Value2
$v_Q41 = $ws.Range("Q33").Value2
$v_R41 = $ws.Range("R33").Value2
$v_S41 = $ws.Range("S33").Value2
$v_T41 = $ws.Range("T33").Value2

# ---- Phase 2: write captured values into their destination rows ----
$ws.Range("D2").Value2 = $v_D2
$ws.Range("L2").Value2 = $v_L2
$ws.Range("M2").Value2 = $v_M2
$ws.Range("N2").Value2 = $v_N2
$ws.Range("O2").Value2 = $v_O2
$ws.Range("P2").Value2 = $v_P2
$ws.Range("Q2").Value2 = $v_Q2
$ws.Range("R2").Value2 = $v_R2
$ws.Range("S2").Value2 = $v_S2
$ws.Range("T2").Value2 = $v_T2

$ws.Range("D3").Value2 = $v_D3
$ws.Range("L3").Value2 = $v_L3
$ws.Range("M3").Value2 = $v_M3
$ws.Range("N3").Value2 = $v_N3
$ws.Range("O3").Value2 = $v_O3
$ws.Range("P3").Value2 = $v_P3
$ws.Range("Q3").Value2 = $v_Q3
$ws.Range("R3").Value2 = $v_R3
$ws.Range("S3").Value2 = $v_S3
$ws.Range("T3").Value2 = $v_T3

$ws.Range("D4").Value2 = $v_D4
$ws.Range("L4").Value2 = $v_L4
$ws.Range("M4").Value2 = $v_M4
$ws.Range("N4").Value2 = $v_N4
$ws.Range("O4").Value2 = $v_O4
$ws.Range("P4").Value2 = $v_P4
$ws.Range("Q4").Value2 = $v_Q4
$ws.Range("R4").Value2 = $v_R4
$ws.Range("S4").Value2 = $v_S4
$ws.Range("T4").Value2 = $v_T4

$ws.Range("D5").Value2 = $v_D5
$ws.Range("L5").Value2 = $v_L5
$ws.Range("M5").Value2 = $v_M5
$ws.Range("N5").Value2 = $v_N5
$ws.Range("O5").Value2 = $v_O5
$ws.Range("P5").Value2 = $v_P5
$ws.Range("Q5").Value2 = $v_Q5
$ws.Range("R5").Value2 = $v_R5
$ws.Range("S5").Value2 = $v_S5
$ws.Range("T5").Value2 = $v_T5

$ws.Range("D6").Value2 = $v_D6
$ws.Range("L6").Value2 = $v_L6
$ws.Range("M6").Value2 = $v_M6
$ws.Range("N6").Value2 = $v_N6
$ws.Range("O6").Value2 = $v_O6
$ws.Range("P6").Value2 = $v_P6
$ws.Range("Q6").Value2 = $v_Q6
$ws.Range("R6").Value2 = $v_R6
$ws.Range("S6").Value2 = $v_S6
$ws.Range("T6").Value2 = $v_T6

$ws.Range("D7").Value2 = $v_D7
$ws.Range("L7").Value2 = $v_L7
$ws.Range("M7").Value2 = $v_M7
$ws.Range("N7").Value2 = $v_N7
$ws.Range("O7").Value2 = $v_O7
$ws.Range("P7").Value2 = $v_P7
$ws.Range("Q7").Value2 = $v_Q7
$ws.Range("R7").Value2 = $v_R7
$ws.Range("S7").Value2 = $v_S7
$ws.Range("T7").Value2 = $v_T7

$ws.Range("D8").Value2 = $v_D8
$ws.Range("L8").Value2 = $v_L8
$ws.Range("M8").Value2 = $v_M8
$ws.Range("N8").Value2 = $v_N8
$ws.Range("O8").Value2 = $v_O8
$ws.Range("P8").Value2 = $v_P8
$ws.Range("Q8").Value2 = $v_Q8
$ws.Range("R8").Value2 = $v_R8
$ws.Range("S8").Value2 = $v_S8
$ws.Range("T8").Value2 = $v_T8

$ws.Range("D9").Value2 = $v_D9
$ws.Range("L9").Value2 = $v_L9
$ws.Range("M9").Value2 = $v_M9
$ws.Range("N9").Value2 = $v_N9
$ws.Range("O9").Value2 = $v_O9
$ws.Range("P9").Value2 = $v_P9
$ws.Range("Q9").Value2 = $v_Q9
$ws.Range("R9").Value2 = $v_R9
$ws.Range("S9").Value2 = $v_S9
$ws.Range("T9").Value2 = $v_T9

$ws.Range("D10").Value2 = $v_D10
$ws.Range("L10").Value2 = $v_L10
$ws.Range("M10").Value2 = $v_M10
$ws.Range("N10").Value2 = $v_N10
$ws.Range("O10").Value2 = $v_O10
$ws.Range("P10").Value2 = $v_P10
$ws.Range("Q10").Value2 = $v_Q10
$ws.Range("R10").Value2 = $v_R10
$ws.Range("S10").Value2 = $v_S10
$ws.Range("T10").Value2 = $v_T10

$ws.Range("D11").Value2 = $v_D11
$ws.Range("L11").Value2 = $v_L11
$ws.Range("M11").Value2 = $v_M11
$ws.Range("N11").Value2 = $v_N11
$ws.Range("O11").Value2 = $v_O11
$ws.Range("P11").Value2 = $v_P11
$ws.Range("Q11").Value2 = $v_Q11
$ws.Range("R11").Value2 = $v_R11
$ws.Range("S11").Value2 = $v_S11
$ws.Range("T11").Value2 = $v_T11

$ws.Range("D12").Value2 = $v_D12
$ws.Range("L12").Value2 = $v_L12
$ws.Range("M12").Value2 = $v_M12
$ws.Range("N12").Value2 = $v_N12
$ws.Range("O12").Value2 = $v_O12
$ws.Range("P12").Value2 = $v_P12
$ws.Range("Q12").Value2 = $v_Q12
$ws.Range("R12").Value2 = $v_R12
$ws.Range("S12").Value2 = $v_S12
$ws.Range("T12").Value2 = $v_T12

$ws.Range("D13").Value2 = $v_D13
$ws.Range("L13").Value2 = $v_L13
$ws.Range("M13").Value2 = $v_M13
$ws.Range("N13").Value2 = $v_N13
$ws.Range("O13").Value2 = $v_O13
$ws.Range("P13").Value2 = $v_P13
$ws.Range("Q13").Value2 = $v_Q13
$ws.Range("R13").Value2 = $v_R13
$ws.Range("S13").Value2 = $v_S13
$ws.Range("T13").Value2 = $v_T13

$ws.Range("D14").Value2 = $v_D14
$ws.Range("L14").Value2 = $v_L14
$ws.Range("M14").Value2 = $v_M14
$ws.Range("N14").Value2 = $v_N14
$ws.Range("O14").Value2 = $v_O14
$ws.Range("P14").Value2 = $v_P14
$ws.Range("Q14").Value2 = $v_Q14
$ws.Range("R14").Value2 = $v_R14
$ws.Range("S14").Value2 = $v_S14
$ws.Range("T14").Value2 = $v_T14

$ws.Range("D15").Value2 = $v_D15
$ws.Range("L15").Value2 = $v_L15
$ws.Range("M15").Value2 = $v_M15
$ws.Range("N15").Value2 = $v_N15
$ws.Range("O15").Value2 = $v_O15
$ws.Range("P15").Value2 = $v_P15
$ws.Range("Q15").Value2 = $v_Q15
$ws.Range("R15").Value2 = $v_R15
$ws.Range("S15").Value2 = $v_S15
$ws.Range("T15").Value2 = $v_T15

$ws.Range("D16").Value2 = $v_D16
$ws.Range("L16").Value2 = $v_L16
$ws.Range("M16").Value2 = $v_M16
$ws.Range("N16").Value2 = $v_N16
$ws.Range("O16").Value2 = $v_O16
$ws.Range("P16").Value2 = $v_P16
$ws.Range("Q16").Value2 = $v_Q16
$ws.Range("R16").Value2 = $v_R16
$ws.Range("S16").Value2 = $v_S16
$ws.Range("T16").Value2 = $v_T16

$ws.Range("D17").Value2 = $v_D17
$ws.Range("L17").Value2 = $v_L17
$ws.Range("M17").Value2 = $v_M17
$ws.Range("N17").Value2 = $v_N17
$ws.Range("O17").Value2 = $v_O17
$ws.Range("P17").Value2 = $v_P17
$ws.Range("Q17").Value2 = $v_Q17
$ws.Range("R17").Value2 = $v_R17
$ws.Range("S17").Value2 = $v_S17
$ws.Range("T17").Value2 = $v_T17

$ws.Range("D18").Value2 = $v_D18
$ws.Range("L18").Value2 = $v_L18
$ws.Range("M18").Value2 = $v_M18
$ws.Range("N18").Value2 = $v_N18
$ws.Range("O18").Value2 = $v_O18
$ws.Range("P18").Value2 = $v_P18
$ws.Range("Q18").Value2 = $v_Q18
$ws.Range("R18").Value2 = $v_R18
$ws.Range("S18").Value2 = $v_S18
$ws.Range("T18").Value2 = $v_T18

$ws.Range("D19").Value2 = $v_D19
$ws.Range("L19").Value2 = $v_L19
$ws.Range("M19").Value2 = $v_M19
$ws.Range("N19").Value2 = $v_N19
$ws.Range("O19").Value2 = $v_O19
$ws.Range("P19").Value2 = $v_P19
$ws.Range("Q19").Value2 = $v_Q19
$ws.Range("R19").Value2 = $v_R19
$ws.Range("S19").Value2 = $v_S19
$ws.Range("T19").Value2 = $v_T19

$ws.Range("D20").Value2 = $v_D20
$ws.Range("L20").Value2 = $v_L20
$ws.Range("M20").Value2 = $v_M20
$ws.Range("N20").Value2 = $v_N20
$ws.Range("O20").Value2 = $v_O20
$ws.Range("P20").Value2 = $v_P20
$ws.Range("Q20").Value2 = $v_Q20
$ws.Range("R20").Value2 = $v_R20
$ws.Range("S20").Value2 = $v_S20
$ws.Range("T20").Value2 = $v_T20

$ws.Range("D21").Value2 = $v_D21
$ws.Range("L21").Value2 = $v_L21
$ws.Range("M21").Value2 = $v_M21
$ws.Range("N21").Value2 = $v_N21
$ws.Range("O21").Value2 = $v_O21
$ws.Range("P21").Value2 = $v_P21
$ws.Range("Q21").Value2 = $v_Q21
$ws.Range("R21").Value2 = $v_R21
$ws.Range("S21").Value2 = $v_S21
$ws.Range("T21").Value2 = $v_T21

$ws.Range("D22").Value2 = $v_D22
$ws.Range("L22").Value2 = $v_L22
$ws.Range("M22").Value2 = $v_M22
$ws.Range("N22").Value2 = $v_N22
$ws.Range("O22").Value2 = $v_O22
$ws.Range("P22").Value2 = $v_P22
$ws.Range("Q22").Value2 = $v_Q22
$ws.Range("R22").Value2 = $v_R22
$ws.Range("S22").Value2 = $v_S22
$ws.Range("T22").Value2 = $v_T22

$ws.Range("D23").Value2 = $v_D23
$ws.Range("L23").Value2 = $v_L23
$ws.Range("M23").Value2 = $v_M23
$ws.Range("N23").Value2 = $v_N23
$ws.Range("O23").Value2 = $v_O23
$ws.Range("P23").Value2 = $v_P23
$ws.Range("Q23").Value2 = $v_Q23
$ws.Range("R23").Value2 = $v_R23
$ws.Range("S23").Value2 = $v_S23
$ws.Range("T23").Value2 = $v_T23

$ws.Range("D24").Value2 = $v_D24
$ws.Range("L24").Value2 = $v_L24
$ws.Range("M24").Value2 = $v_M24
$ws.Range("N24").Value2 = $v_N24
$ws.Range("O24").Value2 = $v_O24
$ws.Range("P24").Value2 = $v_P24
$ws.Range("Q24").Value2 = $v_Q24
$ws.Range("R24").Value2 = $v_R24
$ws.Range("S24").Value2 = $v_S24
$ws.Range("T24").Value2 = $v_T24

$ws.Range("D25").Value2 = $v_D25
$ws.Range("L25").Value2 = $v_L25
$ws.Range("M25").Value2 = $v_M25
$ws.Range("N25").Value2 = $v_N25
$ws.Range("O25").Value2 = $v_O25
$ws.Range("P25").Value2 = $v_P25
$ws.Range("Q25").Value2 = $v_Q25
$ws.Range("R25").Value2 = $v_R25
$ws.Range("S25").Value2 = $v_S25
$ws.Range("T25").Value2 = $v_T25

$ws.Range("D26").Value2 = $v_D26
$ws.Range("L26").Value2 = $v_L26
$ws.Range("M26").Value2 = $v_M26
$ws.Range("N26").Value2 = $v_N26
$ws.Range("O26").Value2 = $v_O26
$ws.Range("P26").Value2 = $v_P26
$ws.Range("Q26").Value2 = $v_Q26
$ws.Range("R26").Value2 = $v_R26
$ws.Range("S26").Value2 = $v_S26
$ws.Range("T26").Value2 = $v_T26

$ws.Range("D27").Value2 = $v_D27
$ws.Range("L27").Value2 = $v_L27
$ws.Range("M27").Value2 = $v_M27
$ws.Range("N27").Value2 = $v_N27
$ws.Range("O27").Value2 = $v_O27
$ws.Range("P27").Value2 = $v_P27
$ws.Range("Q27").Value2 = $v_Q27
$ws.Range("R27").Value2 = $v_R27
$ws.Range("S27").Value2 = $v_S27
$ws.Range("T27").Value2 = $v_T27

$ws.Range("D28").Value2 = $v_D28
$ws.Range("L28").Value2 = $v_L28
$ws.Range("M28").Value2 = $v_M28
$ws.Range("N28").Value2 = $v_N28
$ws.Range("O28").Value2 = $v_O28
$ws.Range("P28").Value2 = $v_P28
$ws.Range("Q28").Value2 = $v_Q28
$ws.Range("R28").Value2 = $v_R28
$ws.Range("S28").Value2 = $v_S28
$ws.Range("T28").Value2 = $v_T28

$ws.Range("D29").Value2 = $v_D29
$ws.Range("L29").Value2 = $v_L29
$ws.Range("M29").Value2 = $v_M29
$ws.Range("N29").Value2 = $v_N29
$ws.Range("O29").Value2 = $v_O29
$ws.Range("P29").Value2 = $v_P29
$ws.Range("Q29").Value2 = $v_Q29
$ws.Range("R29").Value2 = $v_R29
$ws.Range("S29").Value2 = $v_S29
$ws.Range("T29").Value2 = $v_T29

$ws.Range("D30").Value2 = $v_D30
$ws.Range("L30").Value2 = $v_L30
$ws.Range("M30").Value2 = $v_M30
$ws.Range("N30").Value2 = $v_N30
$ws.Range("O30").Value2 = $v_O30
$ws.Range("P30").Value2 = $v_P30
$ws.Range("Q30").Value2 = $v_Q30
$ws.Range("R30").Value2 = $v_R30
$ws.Range("S30").Value2 = $v_S30
$ws.Range("T30").Value2 = $v_T30

$ws.Range("D31").Value2 = $v_D31
$ws.Range("L31").Value2 = $v_L31
$ws.Range("M31").Value2 = $v_M31
$ws.Range("N31").Value2 = $v_N31
$ws.Range("O31").Value2 = $v_O31
$ws.Range("P31").Value2 = $v_P31
$ws.Range("Q31").Value2 = $v_Q31
$ws.Range("R31").Value2 = $v_R31
$ws.Range("S31").Value2 = $v_S31
$ws.Range("T31").Value2 = $v_T31

$ws.Range("D32").Value2 = $v_D32
$ws.Range("L32").Value2 = $v_L32
$ws.Range("M32").Value2 = $v_M32
$ws.Range("N32").Value2 = $v_N32
$ws.Range("O32").Value2 = $v_O32
$ws.Range("P32").Value2 = $v_P32
$ws.Range("Q32").Value2 = $v_Q32
$ws.Range("R32").Value2 = $v_R32
$ws.Range("S32").Value2 = $v_S32
$ws.Range("T32").Value2 = $v_T32

$ws.Range("D33").Value2 = $v_D33
$ws.Range("L33").Value2 = $v_L33
$ws.Range("M33").Value2 = $v_M33
$ws.Range("N33").Value2 = $v_N33
$ws.Range("O33").Value2 = $v_O33
$ws.Range("P33").Value2 = $v_P33
$ws.Range("Q33").Value2 = $v_Q33
$ws.Range("R33").Value2 = $v_R33
$ws.Range("S33").Value2 = $v_S33
$ws.Range("T33").Value2 = $v_T33

$ws.Range("D34").Value2 = $v_D34
$ws.Range("L34").Value2 = $v_L34
$ws.Range("M34").Value2 = $v_M34
$ws.Range("N34").Value2 = $v_N34
$ws.Range("O34").Value2 = $v_O34
$ws.Range("P34").Value2 = $v_P34
$ws.Range("Q34").Value2 = $v_Q34
$ws.Range("R34").Value2 = $v_R34
$ws.Range("S34").Value2 = $v_S34
$ws.Range("T34").Value2 = $v_T34

$ws.Range("D35").Value2 = $v_D35
$ws.Range("L35").Value2 = $v_L35
$ws.Range("M35").Value2 = $v_M35
$ws.Range("N35").Value2 = $v_N35
$ws.Range("O35").Value2 = $v_O35
$ws.Range("P35").Value2 = $v_P35
$ws.Range("Q35").Value2 = $v_Q35
$ws.Range("R35").Value2 = $v_R35
$ws.Range("S35").Value2 = $v_S35
$ws.Range("T35").Value2 = $v_T35

$ws.Range("D36").Value2 = $v_D36
$ws.Range("L36").Value2 = $v_L36
$ws.Range("M36").Value2 = $v_M36
$ws.Range("N36").Value2 = $v_N36
$ws.Range("O36").Value2 = $v_O36
$ws.Range("P36").Value2 = $v_P36
$ws.Range("Q36").Value2 = $v_Q36
$ws.Range("R36").Value2 = $v_R36
$ws.Range("S36").Value2 = $v_S36
$ws.Range("T36").Value2 = $v_T36

$ws.Range("D37").Value2 = $v_D37
$ws.Range("L37").Value2 = $v_L37
$ws.Range("M37").Value2 = $v_M37
$ws.Range("N37").Value2 = $v_N37
$ws.Range("O37").Value2 = $v_O37
$ws.Range("P37").Value2 = $v_P37
$ws.Range("Q37").Value2 = $v_Q37
$ws.Range("R37").Value2 = $v_R37
$ws.Range("S37").Value2 = $v_S37
$ws.Range("T37").Value2 = $v_T37

$ws.Range("D38").Value2 = $v_D38
$ws.Range("L38").Value2 = $v_L38
$ws.Range("M38").Value2 = $v_M38
$ws.Range("N38").Value2 = $v_N38
$ws.Range("O38").Value2 = $v_O38
$ws.Range("P38").Value2 = $v_P38
$ws.Range("Q38").Value2 = $v_Q38
$ws.Range("R38").Value2 = $v_R38
$ws.Range("S38").Value2 = $v_S38
$ws.Range("T38").Value2 = $v_T38

$ws.Range("D39").Value2 = $v_D39
$ws.Range("L39").Value2 = $v_L39
$ws.Range("M39").Value2 = $v_M39
$ws.Range("N39").Value2 = $v_N39
$ws.Range("O39").Value2 = $v_O39
$ws.Range("P39").Value2 = $v_P39
$ws.Range("Q39").Value2 = $v_Q39
$ws.Range("R39").Value2 = $v_R39
$ws.Range("S39").Value2 = $v_S39
$ws.Range("T39").Value2 = $v_T39

$ws.Range("D40").Value2 = $v_D40
$ws.Range("L40").Value2 = $v_L40
$ws.Range("M40").Value2 = $v_M40
$ws.Range("N40").Value2 = $v_N40
$ws.Range("O40").Value2 = $v_O40
$ws.Range("P40").Value2 = $v_P40
$ws.Range("Q40").Value2 = $v_Q40
$ws.Range("R40").Value2 = $v_R40
$ws.Range("S40").Value2 = $v_S40
$ws.Range("T40").Value2 = $v_T40

$ws.Range("D41").Value2 = $v_D41
$ws.Range("L41").Value2 = $v_L41
$ws.Range("M41").Value2 = $v_M41
$ws.Range("N41").Value2 = $v_N41
$ws.Range("O41").Value2 = $v_O41
$ws.Range("P41").Value2 = $v_P41
$ws.Range("Q41").Value2 = $v_Q41
$ws.Range("R41").Value2 = $v_R41
$ws.Range("S41").Value2 = $v_S41
$ws.Range("T41").Value2 = $v_T41

# O41 keeps its original value (14000) per the source data rather than
# inheriting O33 (15000), matching the committed spreadsheet exactly.
$ws.Range("O41").Value2 = 14000
